$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.2038054651530871
$ws.Range("J2").Value = 0.2038054651530872
$ws.Range("M2").Value = 0.52293
$ws.Range("N2").Value = 1.56879
$ws.Range("O2").Value = 0.3445930095145882
$ws.Range("P2").Value = 0.3445930095145882
$ws.Range("Q2").Value = 0.1261429177
$ws.Range("R2").Value = 1.1352862593
$ws.Range("S2").Value = 0.07022993859262282
$ws.Range("T2").Value = 0.07022993859262283

# Row 3
$ws.Range("I3").Value = 0.2038054651530871
$ws.Range("J3").Value = 0.2038054651530872
$ws.Range("O3").Value = 0.4314455865542851
$ws.Range("P3").Value = 0.4314455865542852
$ws.Range("S3").Value = 0.0879309684559426
$ws.Range("T3").Value = 0.08793096845594262

# Row 4
$ws.Range("I4").Value = 0.2038054651530871
$ws.Range("J4").Value = 0.2038054651530872
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2505543333333333
$ws.Range("N4").Value = 0.751663
$ws.Range("O4").Value = 0.1651067480738428
$ws.Range("P4").Value = 0.1651067480738428
$ws.Range("Q4").Value = 0.06043955146777778
$ws.Range("R4").Value = 0.5439559632100001
$ws.Range("S4").Value = 0.03364965759110311
$ws.Range("T4").Value = 0.03364965759110312

# Row 5
$ws.Range("I5").Value = 0.2038054651530871
$ws.Range("J5").Value = 0.2038054651530872
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08931366666666667
$ws.Range("N5").Value = 0.267941
$ws.Range("O5").Value = 0.05885465585728381
$ws.Range("P5").Value = 0.05885465585728381
$ws.Range("Q5").Value = 0.02154454038555556
$ws.Range("R5").Value = 0.19390086347
$ws.Range("S5").Value = 0.01199490051341859
$ws.Range("T5").Value = 0.01199490051341859

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.133983
$ws.Range("H6").Value = 0.401949
$ws.Range("I6").Value = 0.1131999432238703
$ws.Range("J6").Value = 0.1131999432238703
$ws.Range("M6").Value = 0.52293
$ws.Range("N6").Value = 1.56879
$ws.Range("O6").Value = 0.3445930095145882
$ws.Range("P6").Value = 0.3445930095145882
$ws.Range("Q6").Value = 0.07006373018999999
$ws.Range("R6").Value = 0.63057357171
$ws.Range("S6").Value = 0.03900790911239398
$ws.Range("T6").Value = 0.03900790911239398

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.133983
$ws.Range("H7").Value = 0.401949
$ws.Range("I7").Value = 0.1131999432238703
$ws.Range("J7").Value = 0.1131999432238703
$ws.Range("O7").Value = 0.4314455865542851
$ws.Range("P7").Value = 0.4314455865542852
$ws.Range("Q7").Value = 0.087722868234
$ws.Range("R7").Value = 0.789505814106
$ws.Range("S7").Value = 0.04883961590213449
$ws.Range("T7").Value = 0.0488396159021345

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.133983
$ws.Range("H8").Value = 0.401949
$ws.Range("I8").Value = 0.1131999432238703
$ws.Range("J8").Value = 0.1131999432238703
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2505543333333333
$ws.Range("N8").Value = 0.751663
$ws.Range("O8").Value = 0.1651067480738428
$ws.Range("P8").Value = 0.1651067480738428
$ws.Range("Q8").Value = 0.03357002124299999
$ws.Range("R8").Value = 0.302130191187
$ws.Range("S8").Value = 0.01869007450783686
$ws.Range("T8").Value = 0.01869007450783686

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.133983
$ws.Range("H9").Value = 0.401949
$ws.Range("I9").Value = 0.1131999432238703
$ws.Range("J9").Value = 0.1131999432238703
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.08931366666666667
$ws.Range("N9").Value = 0.267941
$ws.Range("O9").Value = 0.05885465585728381
$ws.Range("P9").Value = 0.05885465585728381
$ws.Range("Q9").Value = 0.011966513001
$ws.Range("R9").Value = 0.107698617009
$ws.Range("S9").Value = 0.006662343701504952
$ws.Range("T9").Value = 0.006662343701504953

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8083896666666667
$ws.Range("H10").Value = 2.425169
$ws.Range("I10").Value = 0.6829945916230425
$ws.Range("J10").Value = 0.6829945916230425
$ws.Range("M10").Value = 0.52293
$ws.Range("N10").Value = 1.56879
$ws.Range("O10").Value = 0.3445930095145882
$ws.Range("P10").Value = 0.3445930095145882
$ws.Range("Q10").Value = 0.42273120839
$ws.Range("R10").Value = 3.80458087551
$ws.Range("S10").Value = 0.2353551618095714
$ws.Range("T10").Value = 0.2353551618095714

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8083896666666667
$ws.Range("H11").Value = 2.425169
$ws.Range("I11").Value = 0.6829945916230425
$ws.Range("J11").Value = 0.6829945916230425
$ws.Range("O11").Value = 0.4314455865542851
$ws.Range("P11").Value = 0.4314455865542852
$ws.Range("Q11").Value = 0.5292780443095556
$ws.Range("R11").Value = 4.763502398786
$ws.Range("S11").Value = 0.294675002196208
$ws.Range("T11").Value = 0.2946750021962081

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8083896666666667
$ws.Range("H12").Value = 2.425169
$ws.Range("I12").Value = 0.6829945916230425
$ws.Range("J12").Value = 0.6829945916230425
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2505543333333333
$ws.Range("N12").Value = 0.751663
$ws.Range("O12").Value = 0.1651067480738428
$ws.Range("P12").Value = 0.1651067480738428
$ws.Range("Q12").Value = 0.2025455340052222
$ws.Range("R12").Value = 1.822909806047
$ws.Range("S12").Value = 0.1127670159749028
$ws.Range("T12").Value = 0.1127670159749028

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8083896666666667
$ws.Range("H13").Value = 2.425169
$ws.Range("I13").Value = 0.6829945916230425
$ws.Range("J13").Value = 0.6829945916230425
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.08931366666666667
$ws.Range("N13").Value = 0.267941
$ws.Range("O13").Value = 0.05885465585728381
$ws.Range("P13").Value = 0.05885465585728381
$ws.Range("Q13").Value = 0.07220024522544445
$ws.Range("R13").Value = 0.6498022070289999
$ws.Range("S13").Value = 0.04019741164236026
$ws.Range("T13").Value = 0.04019741164236026
